$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "BOM position" mix-up: row 9 (F9/G9) used to hold the
# mistyped "603-RT0805FRE0747K" mouser part (missing trailing L) together
# with a stale query string, while row 10 (F10/G10) held the correct
# "603-RT0805FRE10750RL" part. Correct row 9 to the properly-typed part
# number and its current Mouser link.
$ws.Range("F9").Value = "603-RT0805FRE0747KL"
$ws.Range("G9").Value = "https://hr.mouser.com/ProductDetail/Yageo/RT0805FRE0747KL?qs=%2Fha2pyFaduhFC1wu9jyUADMq3OpHh%252BVcLIGGOAj8K%2FkOFC1JonHvZN4lfpgCCLyX"

# Row 10 keeps the same text/link it always had; rewriting it (with the
# same values) keeps the two shared-string pairs in the same relative
# order the corrected workbook ends up with.
$ws.Range("F10").Value = "603-RT0805FRE10750RL"
$ws.Range("G10").Value = "https://hr.mouser.com/ProductDetail/Yageo/RT0805FRE10750RL?qs=sGAEpiMZZMvdGkrng054tz3%252BeNFGeSG0b8vprN11rBqrXC%2Fny4PVvQ%3D%3D"

# --- Rebuild the hyperlinks. The underlying Hyperlinks.Delete() call
# clears the whole worksheet collection, so capture the existing targets
# first and re-add every one of them afterwards -- except G9, whose link
# was not recreated for the corrected part (matching the source edit).
$targets = @{
    "G15" = "https://hr.mouser.com/ProductDetail/Micro-Commercial-Components-MCC/SMBJ5356B-TP?qs=%2Fha2pyFadujQlI%252Bovsu2dIVvji%2FULpZcaJMXGhg3YhPF3xL3ld2lcg%3D%3D"
    "G14" = "https://hr.mouser.com/ProductDetail/Analog-Devices/LT1716CS5TRMPBF?qs=%2Fha2pyFadujLvUfqpL8iAa57JmbvGUolaKsSySTrh2UECeK44kz9iw%3D%3D"
    "G13" = "https://hr.mouser.com/ProductDetail/IXYS-Integrated-Circuits/LBA710S?qs=8uBHJDVwVqzgfef1rN5c6w%3D%3D"
    "G12" = "https://hr.mouser.com/ProductDetail/Diodes-Incorporated/DMN10H220L-7?qs=%2Fha2pyFadui6aX5l%2FNbtt80qGfKcPMXPtfI6XoBbWY%252B41%2FteVg33rg%3D%3D"
    "G11" = "https://hr.mouser.com/ProductDetail/Toshiba/CUHS20S30H3F?qs=%2Fha2pyFadugVnzvuqBQjK7mNmFgLGyvkvvKthIvrsa6zTCxP3izpMQ%3D%3D"
    "G4"  = "https://hr.mouser.com/ProductDetail/Diodes-Incorporated/1N4148W-7-F?qs=%2Fha2pyFaduhvdRTUMUAPE5iMJ3chEEQ6fhH0Sc3FvyA%3D"
    "G8"  = "https://hr.mouser.com/ProductDetail/Nichicon/UPW1V472MHD?qs=%2Fha2pyFaduiERwYPL8c3nRShfDk4RNzwZqvNXmDlcV12MjtgXxrA6w%3D%3D"
    "G7"  = "https://hr.mouser.com/ProductDetail/Nichicon/UHD1V331MPD?qs=RhhqrI6N3g%2FaCcBAxpz5YQ%3D%3D"
    "G6"  = "https://hr.mouser.com/ProductDetail/Nichicon/UHE1V332MHD6?qs=%2Fha2pyFaduh%2FVvtnmOexBNcjzEqxV8XETiySkoDb%2Fy5n4ZjrNuYHaA%3D%3D"
    "G3"  = "https://hr.mouser.com/ProductDetail/Yageo/RT0805FRE0710KL?qs=sGAEpiMZZMvdGkrng054t%252BKCHBXLTLydi6xJ2%2FVOu7c%3D"
    "G2"  = "https://hr.mouser.com/ProductDetail/Yageo/RT0805BRD07100KL?qs=sGAEpiMZZMvdGkrng054t%252BKCHBXLTLydbVwVVm%252B5HQE%3D"
    "G5"  = "https://hr.mouser.com/ProductDetail/Yageo/RT0805BRD0730KL?qs=sGAEpiMZZMvdGkrng054t%252BKCHBXLTLydlJaovHNoOLo%3D"
    "G10" = "https://hr.mouser.com/ProductDetail/Yageo/RT0805FRE10750RL?qs=sGAEpiMZZMvdGkrng054tz3%252BeNFGeSG0b8vprN11rBqrXC%2Fny4PVvQ%3D%3D"
}
$order = @("G15", "G14", "G13", "G12", "G11", "G4", "G8", "G7", "G6", "G3", "G2", "G5", "G10")

# Hyperlinks.Add() re-applies the built-in "Hyperlink" font/style to its
# target cell, which would otherwise drift G2:G8/G10:G15 away from the
# plain style index they already carried (s="2"). Snapshot + restore it
# so only the intended G9/G10 content actually changes.
$origStyles = @{}
foreach ($ref in $order) {
    $origStyles[$ref] = $ws.Range($ref).Style
}

$ws.Hyperlinks.Delete()
foreach ($ref in $order) {
    $ws.Hyperlinks.Add($ws.Range($ref), $targets[$ref])
}
foreach ($ref in $order) {
    $ws.Range($ref).Style = $origStyles[$ref]
}

# --- Cosmetic: cursor moved to D23 as part of the edit.
$ws.Range("D23").Select()
